$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = 5
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 5

$ws.Range("F25").Select()
